$wb = $excel.ActiveWorkbook

# --- "Repayment schedule" sheet: insert a new (blank) column before column N ---
$ws = $wb.Worksheets.Item("Repayment schedule")

# Inserting a column shifts existing N:P -> O:Q and leaves the new N column blank,
# inheriting the formatting of the column to its left (M).
$ws.Columns("N").Insert()
$ws.Columns("N").ColumnWidth = $ws.Columns("M").ColumnWidth

# Update the selected cell on this sheet
$ws.Range("S7").Select() | Out-Null

# Make "Repayment schedule" the active sheet/tab (was "Summary" before)
$ws.Activate() | Out-Null
